$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record needs to be inserted at row 121, pushing the
# existing rows 121-130 down to 122-131 (the sheet's last row moves from
# 130 to 131).
$ws.Rows.Item(121).Insert()

# Populate the newly inserted row 121 with the new Papaya price record.
$ws.Range("A121").Value = 10
$ws.Range("B121").Value = "Vega Modelo de Temuco"
$ws.Range("C121").Value = "La Araucanía"
$ws.Range("D121").Value = 45223
$ws.Range("E121").Value = 9
$ws.Range("F121").Value = "Fruta"
$ws.Range("G121").Value = 100108
$ws.Range("H121").Value = "Tropicales y subtropicales"
$ws.Range("I121").Value = 100108004
$ws.Range("J121").Value = "Papaya"
$ws.Range("K121").Value = "Cultivar IV Región"
$ws.Range("L121").Value = "Primera"
$ws.Range("M121").Value = 60
$ws.Range("N121").Value = 24000
$ws.Range("O121").Value = 24000
$ws.Range("P121").Value = 24000
$ws.Range("Q121").Value = "$/bandeja 10 kilos"
$ws.Range("R121").Value = "Provincia del Elquí"
$ws.Range("S121").Value = 2400
$ws.Range("T121").Value = 10
